$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Elizabeth J.%Anderson%xref no email%1,    Kacey C.%Ernst%xref no email%1,    Francisco Fernando%Martins%xref no email%1,    Cicera da Silva%Martins%xref no email%1,    Mary P.%Koss%xref no email%1]"
$ws.Range("E3").Value = "[Ana Rosa%Linde-Arias%linde14@yahoo.com%1,            Maria%Roura%NULL%1,            Eduardo%Siqueira%NULL%1]"
$ws.Range("E4").Value = "[Maria%Belizan%mbelizan@iecs.org.ar%1,            Edna%Maradiaga%edjamar3006@yahoo.com%1,            Javier%Roberti%jroberti@iecs.org.ar%1,            Maricela%Casco-Aguilar%marykasco@yahoo.com%1,            Alison F.%Ortez%alison_fabiola@yahoo.es%1,            Juan C.%Avila-Flores%javilaflores3@gmail.com%1,            Gloria%González%marilyntoin@yahoo.com%1,            Carolina%Bustillo%mcbu1502@yahoo.com%0,            Alejandra%Calderón%lilianalecalderon@gmail.com%1,            Harry%Bock%hbockme@hotmail.com%1,            María L.%Cafferata%NULL%1,            Adriano B.%Tavares%adriano.b.tavares@gmail.com%1,            Jackeline%Alger%jackelinealger@gmail.com%1,            Moazzam%Ali%alimoa@who.int%1]"
$ws.Range("D5").Value = "`nThe primary purpose of this qualitative study was to explore the knowledge, attitudes, and behaviors of low-income women considered high priority for receiving the novel influenza A (H1N1) vaccine to improve communication in emergency preparedness and response.`n Researchers sought to identify the factors that affect this high priority population’s ability to successfully comply with vaccination recommendations.`n By utilizing an existing communication framework through the special supplemental nutrition program for women, infants, and children (WIC) they were able to document the systems and infrastructure needed to foster constructive responses in a sustainable manner in the future.`n Six focus group discussions with WIC clients (n = 56) and 10 individual interviews with staff members were conducted at two WIC clinics in Georgia (1 urban and 1 rural).`n Data were collected after the 2009–2010 influenza season and analyzed using thematic analysis.`n Knowledge and attitudes regarding H1N1 differed among participants with regard to perceived severity and perceived risk of influenza illness.`n Participants identified several barriers and motivators to receiving the vaccination, as well as information needs, sources, and information-seeking behaviors.`n Similarities emerged among both WIC clients and staff members regarding impressions of H1N1 and the vaccine’s use, suggesting that while the information may be provided, it is not effectively understood or accepted.`n Comprehensive education, policy and planning development regarding pandemic influenza and vaccine acceptance among low-income women is necessary, including improvements in risk communication messages and identifying effective methods to disseminate trusted information to these high priority groups.`n"
$ws.Range("E5").Value = "[ Catherine A.%Boyd%null%0,  Julie A.%Gazmararian%null%1,  Winifred Wilkins%Thompson%null%1]"
$ws.Range("E6").Value = "[Benner%P.%coreGivesNoEmail%1,           Bohme%G.%coreGivesNoEmail%1,           Chiang%H-H%coreGivesNoEmail%4,           Chiang%H-H%coreGivesNoEmail%0,           Chiang%H-H%coreGivesNoEmail%0,           Chiang%H-H%coreGivesNoEmail%0,           Donohoe%J.%coreGivesNoEmail%1,           Hsien-Hsien%Chiang%coreGivesNoEmail%1,           I-Ling%Sue%coreGivesNoEmail%1,           International%Council of Nurses%coreGivesNoEmail%1,           Mei-Bih%Chen%coreGivesNoEmail%1,           Morgan%DL%coreGivesNoEmail%1,           Packer%MJ%coreGivesNoEmail%1]"
$ws.Range("E7").Value = "[Chee%YO%coreGivesNoEmail%1,           Dodgson%JE%coreGivesNoEmail%1,           Tarrant%M%coreGivesNoEmail%1,           Watkins%A%coreGivesNoEmail%1]"
$ws.Range("E8").Value = "[ Elin%Erland%null%1,             Bente%Dahl%null%1,           Elin%Erland%null%1,           Bente%Dahl%null%1]"
$ws.Range("E9").Value = "[Davis%Mark%coreGivesNoEmail%1,           Flowers%Paul%coreGivesNoEmail%1,           Lohm%Davina%coreGivesNoEmail%1,           Stephenson%Niamh%coreGivesNoEmail%1,           Waller%Emily%coreGivesNoEmail%1]"
$ws.Range("E10").Value = "[Hector M.%Gomez%NULL%1,            Carlos%Mejia Arbelaez%NULL%1,            Jovana A.%Ocampo Cañas%ja.ocampo@uniandes.edu.co%1]"
$ws.Range("E11").Value = "[ Susan%Jones%null%1,             Betty%Sam%null%1,             Florence%Bull%null%1,             Steven Bagie%Pieh%null%1,             Jaki%Lambert%null%1,             Florence%Mgawadere%null%1,             Somasundari%Gopalakrishnan%null%1,             Charles A.%Ameh%null%1,             Nynke%van den Broek%null%1,           Susan%Jones%null%1,           Betty%Sam%null%1,           Florence%Bull%null%1,           Steven Bagie%Pieh%null%1,           Jaki%Lambert%null%1,           Florence%Mgawadere%null%1,           Somasundari%Gopalakrishnan%null%1,           Charles A.%Ameh%null%1,           Nynke%van den Broek%null%1]"
$ws.Range("E12").Value = "[Caroline S.E.%Homer%NULL%1,            Miranda%Davies-Tuck%NULL%1,            Hannah G.%Dahlen%NULL%1,            Vanessa L.%Scarf%NULL%1]"
$ws.Range("E13").Value = "[Ana Rosa%Linde%xref no email%1,    Carlos Eduardo%Siqueira%xref no email%1]"
$ws.Range("D14").Value = "`nThe objective of this study was to explore pregnant and recently pregnant women’s perceptions of influenza vaccine and antivirals during the 2009 H1N1 pandemic.`n We conducted 18 focus groups with pregnant and recently pregnant women in three US cities in September 2009. Participants were segmented into groups by insurance status (no or public insurance vs.`n private insurance), vaccine attitudes (higher vs.`n lower likelihood of acceptance of any vaccines, not only influenza vaccines), and parity (first child vs.`n other children in the home) based on information they provided on the screening questionnaire at the time of recruitment.`n We found that women are not well informed about influenza vaccinations and antiviral medicine and have significant concerns about taking them during pregnancy.`n An interest in their infant’s well-being, however, can be strong motivation to adopt preventive recommendations, including vaccination.`n A woman’s health care provider is a highly trusted source of information about the 2009 H1N1. Pregnant women have unique communication needs for influenza.`n Messages directing pregnant women to adopt public health recommendations, particularly for vaccination or prophylactic medication should include a detailed description of the benefits or lack of risk to the fetus and the safety of breastfeeding.`n Additionally, messages should recognize that pregnant women are taught to be selective about taking medication and provide a clear rationale for why the medicine or vaccine is necessary.`n"
$ws.Range("E14").Value = "[ Molly M.%Lynch%null%0,  Elizabeth W.%Mitchell%null%1,  Jennifer L.%Williams%null%1,  Kelly%Brumbaugh%null%1,  Michelle%Jones-Bell%null%1,  Debra E.%Pinkney%null%1,  Christine M.%Layton%null%1,  Patricia W.%Mersereau%null%1,  Juliette S.%Kendrick%null%1,  Paula Eguino%Medina%null%1,  Lucia Rojas%Smith%null%1]"
$ws.Range("E15").Value = "[Niuniu%Sun%NULL%1,            Luoqun%Wei%NULL%1,            Suling%Shi%NULL%1,            Dandan%Jiao%NULL%1,            Runluo%Song%NULL%1,            Lili%Ma%NULL%1,            Hongwei%Wang%NULL%1,            Chao%Wang%NULL%1,            Zhaoguo%Wang%NULL%1,            Yanli%You%NULL%1,            Shuhua%Liu%NULL%1,            Hongyun%Wang%NULL%1]"
$ws.Range("E16").Value = "[Bushra%Sabri%xref no email%1,    Maria%Hartley%xref no email%1,    Jyoti%Saha%xref no email%1,    Sarah%Murray%xref no email%1,    Nancy%Glass%xref no email%1,    Jacquelyn C.%Campbell%xref no email%1]"
$ws.Range("E17").Value = "[Veronika%Tirado%Veronika.Tirado@ki.se%1,            Santiago A.%Morales Mesa%NULL%1,            John%Kinsman%NULL%1,            Anna Mia%Ekström%NULL%1,            Berta N.%Restrepo Jaramillo%NULL%1]"
$ws.Range("E18").Value = "[Simon N%Williams%NULL%1,            Christopher J%Armitage%NULL%2,            Christopher J%Armitage%NULL%0,            Tova%Tampe%NULL%1,            Kimberly%Dienes%NULL%1]"
